$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999904767024664
$ws.Range("E2").Value = 0.9999904767024664

# Row 3
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = 0.02106879882139961
$ws.Range("E3").Value = 0.02106879882139961

# Row 4
$ws.Range("D4").Value = 0.9997110917106176
$ws.Range("E4").Value = 0.9997110917106176

# Row 5
$ws.Range("D5").Value = 0.9997931785426536
$ws.Range("E5").Value = 0.9997931785426536

# Row 6
$ws.Range("D6").Value = 0.9999999988129999
$ws.Range("E6").Value = 0.9999999988129999

# Row 7
$ws.Range("D7").Value = 0.9999999998503226
$ws.Range("E7").Value = [double]"1.496773816000996E-10"

# Row 8
$ws.Range("D8").Value = [double]"2.064997771185897E-05"
$ws.Range("E8").Value = 0.9999793500222881

# Row 9
$ws.Range("D9").Value = 0.9791823476382642
$ws.Range("E9").Value = 0.02081765236173583

# Row 11
$ws.Range("D11").Value = 0.9999999999284439
$ws.Range("E11").Value = [double]"7.155609438314059E-11"
$ws.Range("F11").Value = 5.95767879486084
$ws.Range("G11").Value = 0.5
